# "Fruta / hortaliza, semanal" - re-sequence the weekly price records.
# The sheet lists one row per market day; this edit re-shuffles the
# per-day fields (Fecha/Volumen/Precio minimo/Precio promedio ponderado/
# Precio $/Kg) across the existing rows so each row ends up showing a
# different week's figures while the descriptive columns (market,
# region, category, unit, origin, etc.) stay put.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# row -> row whose D/J/K/M/P values it should adopt
$mapping = @{
    2  = 4
    3  = 11
    4  = 15
    5  = 5
    6  = 10
    7  = 19
    8  = 7
    9  = 13
    10 = 3
    11 = 16
    12 = 12
    13 = 8
    14 = 2
    15 = 6
    16 = 9
    17 = 20
    18 = 18
    19 = 14
    20 = 17
}

$cols = @("D", "J", "K", "M", "P")

# Snapshot the current (pre-edit) values for the columns that move,
# so the re-assignment below always reads the original data even
# after earlier rows have already been overwritten.
$original = @{}
foreach ($row in 2..20) {
    $rowVals = @{}
    foreach ($col in $cols) {
        $rowVals[$col] = $ws.Range("$col$row").Value()
    }
    $original[$row] = $rowVals
}

foreach ($row in 2..20) {
    $srcRow = $mapping[$row]
    $srcVals = $original[$srcRow]
    foreach ($col in $cols) {
        $ws.Range("$col$row").Value = $srcVals[$col]
    }
}
